# Slide 4: "What is a "Trust Score"" - update the stated number of
# statistical risk factors from 42 to 62 in the body placeholder,
# editing only that paragraph's run so the rest of the text box is
# left untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item("Content Placeholder 2")
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$run = $para.Runs(1, 1)
$run.Text = "Statistically calculated using 62 factors across 7 risk categories"
